$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "person 1"
$ws.Range("A3").Value = "person 2"
$ws.Range("A4").Value = "person 3"

$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "job"

$ws.Range("B2").Value = "salesman"
$ws.Range("B3").Value = "accountant"
$ws.Range("B4").Value = "programmer"

$ws.Range("C1").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 3

$ws.Range("C4").Select()
